$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B5").Value = 580.50532224355379
$ws.Range("C5").Value = 9.3068926943918093
$ws.Range("D5").Value = 70.867488554029336

$ws.Range("B7").Value = 1490.6151261202419
$ws.Range("C7").Value = 8.9343719981487766

$ws.Range("B8").Value = 927.44315655741389
$ws.Range("C8").Value = 11.732132683963624
$ws.Range("D8").Value = 162.46166697199496

$ws.Range("B11").Value = 351.06068485152878
$ws.Range("D11").Value = 210.86204638883373

$ws.Range("B12").Value = 432.41164549458517
$ws.Range("D12").Value = 179.94861947503836

$ws.Range("B14").Value = 527.91300437905579
$ws.Range("C14").Value = 6.4847493837403682
$ws.Range("D14").Value = 50.467769472570197

$ws.Range("B17").Value = 521.50110307385091
$ws.Range("C17").Value = 4.1108949928782179
$ws.Range("D17").Value = 168.42384845518143

$ws.Range("B18").Value = 711.99448681809326
$ws.Range("C18").Value = 3.8846492643499761
$ws.Range("D18").Value = 110.24303429960459

$ws.Range("B19").Value = 355.05525296766757
$ws.Range("D19").Value = 239.76070478080715

$ws.Range("B20").Value = 958.71609887123418
$ws.Range("C20").Value = 13.70941693504064
$ws.Range("D20").Value = 71.940203442177179

$ws.Range("B21").Value = 388.13948223176465

$ws.Range("B22").Value = 401.94178429995532
$ws.Range("D22").Value = 215.38428967248501

$ws.Range("B23").Value = 376.73028859011436
$ws.Range("D23").Value = 245.37276990622675

$ws.Range("B24").Value = 350.82560286422847
$ws.Range("D24").Value = 245.81792060202173

$ws.Range("B25").Value = 349.90495929360804
$ws.Range("D25").Value = 212.90806822949787

$ws.Range("B26").Value = 663.27883183403026
$ws.Range("C26").Value = 9.9758826732811254
$ws.Range("D26").Value = 102.83414237765939

$ws.Range("B27").Value = 1021.5450896529127
$ws.Range("C27").Value = 12.358464535243575
$ws.Range("D27").Value = 66.381992731925493

$ws.Range("B28").Value = 1264.0236147957798
$ws.Range("C28").Value = 14.263278763576519
$ws.Range("D28").Value = 106.71955488022718
